$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers ---
$ws.Range("H1").Value = "panel SVG"

# Row 2: F2 becomes "diameter", old "top extra" text moves to G2
$ws.Range("F2").Value = "diameter"
$ws.Range("G2").Value = "top extra"
$ws.Range("H2").Value = "X"
$ws.Range("I2").Value = "Y"
$ws.Range("J2").Value = "diameter"

# --- Row 3 (origin) : move top-extra value from F3 to G3, add H3/I3 ---
$ws.Range("F3").ClearContents()
$ws.Range("G3").Value = 10.16
$ws.Range("H3").Value = 50
$ws.Range("I3").Value = 50

# --- New F column (diameter) values, rows 4-12 ---
$ws.Range("F4").Value = 5
$ws.Range("F5").Value = 9.5
$ws.Range("F6").Value = 7.3
$ws.Range("F7").Value = 7.3
$ws.Range("F8").Value = 7.3
$ws.Range("F9").Value = 7.3
$ws.Range("F10").Value = 6.3
$ws.Range("F11").Value = 6.3
$ws.Range("F12").Value = 6.3

# --- New J column (diameter) values, rows 4-12 ---
$ws.Range("J4").Value = 5
$ws.Range("J5").Value = 9.5
$ws.Range("J6").Value = 7
$ws.Range("J7").Value = 7
$ws.Range("J8").Value = 7
$ws.Range("J9").Value = 7
$ws.Range("J10").Value = 6
$ws.Range("J11").Value = 6
$ws.Range("J12").Value = 6

# --- Update existing E formulas to reference $G$3 instead of $F$3 (rows 4-12) ---
for ($r = 4; $r -le 12; $r++) {
    $ws.Range("E$r").Formula = "=C$r-C`$3+E`$3+`$G`$3"
}

# --- New H/I column formulas, rows 4-12 ---
for ($r = 4; $r -le 12; $r++) {
    $ws.Range("H$r").Formula = "=D$r-D`$3+H`$3-`$J$r/2"
    $ws.Range("I$r").Formula = "=E$r-E`$3+I`$3-`$J$r/2"
}

# --- Update selection to match target workbook state ---
$ws.Range("H12").Select()
